$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.877.00'
$ws.Range("E2").Value = '  -4.93%  '
$ws.Range("D3").Value = '2.207.90'
$ws.Range("E3").Value = '  -6.46%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '''316.13'
$ws.Range("E5").Value = '  +1.38%  '
$ws.Range("D6").Value = '''98.90'
$ws.Range("E6").Value = '  -7.96%  '
$ws.Range("E7").Value = '  -6.95%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("E9").Value = '  -8.24%  '
$ws.Range("D10").Value = '''36.78'
$ws.Range("E10").Value = '  -9.89%  '
$ws.Range("D11").Value = '''54.01'
$ws.Range("E11").Value = '  -2.96%  '
$ws.Range("D12").Value = '''0.0824'
$ws.Range("E12").Value = '  -9.99%  '
$ws.Range("D13").Value = '''7.67'
$ws.Range("E13").Value = '  -9.25%  '
$ws.Range("E14").Value = '  -1.82%  '
$ws.Range("D15").Value = '''0.857'
$ws.Range("E15").Value = '  -11.87%  '
$ws.Range("D16").Value = '2.542.06'
$ws.Range("E16").Value = '  -6.55%  '
$ws.Range("D17").Value = '''14.18'
$ws.Range("E17").Value = '  -6.69%  '
$ws.Range("D18").Value = '2.201.72'
$ws.Range("E18").Value = '  -6.57%  '
$ws.Range("D19").Value = '42.785.14'
$ws.Range("E19").Value = '  -5.13%  '
$ws.Range("D20").Value = '''14.37'
$ws.Range("E20").Value = '  +1.11%  '
$ws.Range("D21").Value = '0.0₃0958'
$ws.Range("E21").Value = '  -9.64%  '
$ws.Range("D22").Value = '''6.40'
$ws.Range("E22").Value = '  -10.79%  '
$ws.Range("D23").Value = '''65.16'
$ws.Range("E23").Value = '  -11.01%  '
$ws.Range("D24").Value = '''3.15'
$ws.Range("E24").Value = '  -10.24%  '
$ws.Range("D25").Value = '''235.51'
$ws.Range("E25").Value = '  -9.01%  '
$ws.Range("D26").Value = '''2.11'
$ws.Range("E26").Value = '  -8.31%  '
$ws.Range("E27").Value = '  -0.37%  '
$ws.Range("E28").Value = '  +1.76%  '
$ws.Range("D29").Value = '''9.96'
$ws.Range("E29").Value = '  -9.58%  '
$ws.Range("E30").Value = '  -4.09%  '
$ws.Range("E31").Value = '  -12.52%  '
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").Value = '''0.0884'
$ws.Range("E32").Value = '  -9.27%  '
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").Value = '''20.49'
$ws.Range("E33").Value = '  -8.11%  '
$ws.Range("D34").Value = '''34.27'
$ws.Range("E34").Value = '  -7.86%  '
$ws.Range("D35").Value = '''154.31'
$ws.Range("E35").Value = '  -7.84%  '
$ws.Range("E36").Value = '  -6.85%  '
$ws.Range("D37").Value = '''3.07'
$ws.Range("E37").Value = '  +7.16%  '
$ws.Range("E38").Value = '  -6.89%  '
$ws.Range("E39").Value = '  +7.21%  '
$ws.Range("D40").Value = '''0.107'
$ws.Range("E40").Value = '  -7.47%  '
$ws.Range("D41").Value = '''4.41'
$ws.Range("E41").Value = '  -5.56%  '
$ws.Range("E42").Value = '  -4.44%  '
$ws.Range("E43").Value = '  -8.07%  '
$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").Value = '''1.00'
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '1.839.01'
$ws.Range("E45").Value = '  -0.03%  '
$ws.Range("D46").Value = '''12.23'
$ws.Range("E46").Value = '  -4.59%  '
$ws.Range("D47").Value = '''87.91'
$ws.Range("E47").Value = '  -11.93%  '
$ws.Range("E48").Value = '  -9.16%  '
$ws.Range("D49").Value = '''5.35'
$ws.Range("E49").Value = '  -6.09%  '
$ws.Range("D50").Value = '''60.61'
$ws.Range("E50").Value = '  -12.82%  '
$ws.Range("D51").Value = '''75.38'
$ws.Range("E51").Value = '  -9.69%  '
